$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated HH tide values (spring-moon date plus next-day lookup now takes
# the max of all returned entries instead of just the first one).
$ws.Range("B3").Value  = 5.768
$ws.Range("B4").Value  = 7.1
$ws.Range("B5").Value  = 6.106
$ws.Range("B6").Value  = 6.762
$ws.Range("B7").Value  = 6.01
$ws.Range("B8").Value  = 6.749
$ws.Range("B9").Value  = 5.732
$ws.Range("B10").Value = 6.342
$ws.Range("B11").Value = 6.112
$ws.Range("B13").Value = 6.985
$ws.Range("B15").Value = 6.588
$ws.Range("B16").Value = 6.283
$ws.Range("B17").Value = 7.149
$ws.Range("B20").Value = 5.089
$ws.Range("B21").Value = 6.352
$ws.Range("B25").Value = 6.739
$ws.Range("B28").Value = 6.25
$ws.Range("B29").Value = 6.351520000000002

# B3 previously had the default (no explicit) style; it now matches the
# "0.000" number-format style used by the rest of the B column.
$ws.Range("B3").NumberFormat = $ws.Range("B4").NumberFormat
